$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 55
$ws.Cells.Item(55, 1).Value = "WGG 02"
$ws.Cells.Item(55, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(55, 3).Value = "13-01-2026"
$ws.Cells.Item(55, 4).Value = 286962
$ws.Cells.Item(55, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(55, 6).Value = 34413429360
$ws.Cells.Item(55, 7).Value = "NEFT"
$ws.Cells.Item(55, 8).Value = "SBIN0003229"
$ws.Cells.Item(55, 9).Value = "AAAFW8862C"
$ws.Cells.Item(55, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(55, 11).Value = ""
$ws.Cells.Item(55, 12).Value = "274e755a-9f8a-4a84-a81a-f51955d61f9e"
$ws.Cells.Item(55, 13).Value = ""
$ws.Cells.Item(55, 14).Value = ""
$ws.Cells.Item(55, 15).Value = ""
$ws.Cells.Item(55, 16).Value = ""
$ws.Cells.Item(55, 17).Value = ""
$ws.Cells.Item(55, 18).Value = ""
$ws.Cells.Item(55, 19).Value = ""
$ws.Cells.Item(55, 20).Value = ""
$ws.Cells.Item(55, 21).Value = "pending"
$ws.Cells.Item(55, 22).Value = 10107
$ws.Cells.Item(55, 23).Value = ""
$ws.Cells.Item(55, 24).Value = "Hyrider Service amount Today 12/01/2026 RPA_UNIQUE_ID : ae4cde02-0c4a-4bb1-99e2-a94529153273"
$ws.Cells.Item(55, 25).Value = "HO"
$ws.Cells.Item(55, 26).Value = "PAYMENT"
$ws.Cells.Item(55, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(55, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(55, 29).Value = 0
$ws.Cells.Item(55, 30).Value = 0
$ws.Cells.Item(55, 31).Value = 0
$ws.Cells.Item(55, 32).Value = ""
$ws.Cells.Item(55, 33).Value = ""
$ws.Cells.Item(55, 34).Value = ""
$ws.Cells.Item(55, 35).Value = ""
$ws.Cells.Item(55, 36).Value = ""
$ws.Cells.Item(55, 37).Value = ""
$ws.Cells.Item(55, 38).Value = ""
$ws.Cells.Item(55, 39).Value = ""
$ws.Cells.Item(55, 40).Value = ""
$ws.Cells.Item(55, 41).Value = ""

# Row 56
$ws.Cells.Item(56, 1).Value = "WGG 02"
$ws.Cells.Item(56, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(56, 3).Value = "13-01-2026"
$ws.Cells.Item(56, 4).Value = 286962
$ws.Cells.Item(56, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(56, 6).Value = 34413429360
$ws.Cells.Item(56, 7).Value = "NEFT"
$ws.Cells.Item(56, 8).Value = "SBIN0003229"
$ws.Cells.Item(56, 9).Value = "AAAFW8862C"
$ws.Cells.Item(56, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(56, 11).Value = ""
$ws.Cells.Item(56, 12).Value = "f43a8e50-3fa5-42c2-a6e1-b73ce012316e"
$ws.Cells.Item(56, 13).Value = ""
$ws.Cells.Item(56, 14).Value = ""
$ws.Cells.Item(56, 15).Value = ""
$ws.Cells.Item(56, 16).Value = ""
$ws.Cells.Item(56, 17).Value = ""
$ws.Cells.Item(56, 18).Value = ""
$ws.Cells.Item(56, 19).Value = ""
$ws.Cells.Item(56, 20).Value = ""
$ws.Cells.Item(56, 21).Value = "pending"
$ws.Cells.Item(56, 22).Value = 7953
$ws.Cells.Item(56, 23).Value = ""
$ws.Cells.Item(56, 24).Value = "Hisham sir flat Electricity bill Due on Jan 19, 2026 RPA_UNIQUE_ID : 0a5dc5d6-14c0-48a4-98c2-a10011a7e3e3"
$ws.Cells.Item(56, 25).Value = "HO"
$ws.Cells.Item(56, 26).Value = "PAYMENT"
$ws.Cells.Item(56, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(56, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(56, 29).Value = 0
$ws.Cells.Item(56, 30).Value = 0
$ws.Cells.Item(56, 31).Value = 0
$ws.Cells.Item(56, 32).Value = ""
$ws.Cells.Item(56, 33).Value = ""
$ws.Cells.Item(56, 34).Value = ""
$ws.Cells.Item(56, 35).Value = ""
$ws.Cells.Item(56, 36).Value = ""
$ws.Cells.Item(56, 37).Value = ""
$ws.Cells.Item(56, 38).Value = ""
$ws.Cells.Item(56, 39).Value = ""
$ws.Cells.Item(56, 40).Value = ""
$ws.Cells.Item(56, 41).Value = ""

# Row 57
$ws.Cells.Item(57, 1).Value = "WGE 97"
$ws.Cells.Item(57, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(57, 3).Value = "13-01-2026"
$ws.Cells.Item(57, 4).Value = 286962
$ws.Cells.Item(57, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(57, 6).Value = 34413429360
$ws.Cells.Item(57, 7).Value = "NEFT"
$ws.Cells.Item(57, 8).Value = "SBIN0003229"
$ws.Cells.Item(57, 9).Value = "AAAFW8862C"
$ws.Cells.Item(57, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(57, 11).Value = "Zillahmol S"
$ws.Cells.Item(57, 12).Value = "6f5f73da-810c-4ac1-b979-8b27c4282191"
$ws.Cells.Item(57, 13).Value = 395502010022654
$ws.Cells.Item(57, 14).Value = "UBIN0826308"
$ws.Cells.Item(57, 15).Value = ""
$ws.Cells.Item(57, 16).Value = ""
$ws.Cells.Item(57, 17).Value = ""
$ws.Cells.Item(57, 18).Value = ""
$ws.Cells.Item(57, 19).Value = ""
$ws.Cells.Item(57, 20).Value = ""
$ws.Cells.Item(57, 21).Value = "pending"
$ws.Cells.Item(57, 22).Value = 22000
$ws.Cells.Item(57, 23).Value = ""
$ws.Cells.Item(57, 24).Value = "2 second hand laptop purchase RPA_UNIQUE_ID : 32571bb6-0e7b-489b-abcd-635729bce130"
$ws.Cells.Item(57, 25).Value = "HO"
$ws.Cells.Item(57, 26).Value = "PAYMENT"
$ws.Cells.Item(57, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(57, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(57, 29).Value = 0
$ws.Cells.Item(57, 30).Value = 0
$ws.Cells.Item(57, 31).Value = 0
$ws.Cells.Item(57, 32).Value = ""
$ws.Cells.Item(57, 33).Value = ""
$ws.Cells.Item(57, 34).Value = ""
$ws.Cells.Item(57, 35).Value = ""
$ws.Cells.Item(57, 36).Value = ""
$ws.Cells.Item(57, 37).Value = ""
$ws.Cells.Item(57, 38).Value = ""
$ws.Cells.Item(57, 39).Value = ""
$ws.Cells.Item(57, 40).Value = ""
$ws.Cells.Item(57, 41).Value = ""

# Row 58
$ws.Cells.Item(58, 1).Value = "WGE 97"
$ws.Cells.Item(58, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(58, 3).Value = "13-01-2026"
$ws.Cells.Item(58, 4).Value = 286962
$ws.Cells.Item(58, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(58, 6).Value = 34413429360
$ws.Cells.Item(58, 7).Value = "NEFT"
$ws.Cells.Item(58, 8).Value = "SBIN0003229"
$ws.Cells.Item(58, 9).Value = "AAAFW8862C"
$ws.Cells.Item(58, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(58, 11).Value = "Zillahmol S"
$ws.Cells.Item(58, 12).Value = "cb6dffe3-6402-433c-8e29-d3be104759dd"
$ws.Cells.Item(58, 13).Value = 395502010022654
$ws.Cells.Item(58, 14).Value = "UBIN0826308"
$ws.Cells.Item(58, 15).Value = ""
$ws.Cells.Item(58, 16).Value = ""
$ws.Cells.Item(58, 17).Value = ""
$ws.Cells.Item(58, 18).Value = ""
$ws.Cells.Item(58, 19).Value = ""
$ws.Cells.Item(58, 20).Value = ""
$ws.Cells.Item(58, 21).Value = "pending"
$ws.Cells.Item(58, 22).Value = 2300
$ws.Cells.Item(58, 23).Value = ""
$ws.Cells.Item(58, 24).Value = "Conference room AC repair RPA_UNIQUE_ID : 3d6438e2-61f9-4948-9f8d-713b84d7ea57"
$ws.Cells.Item(58, 25).Value = "HO"
$ws.Cells.Item(58, 26).Value = "PAYMENT"
$ws.Cells.Item(58, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(58, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(58, 29).Value = 0
$ws.Cells.Item(58, 30).Value = 0
$ws.Cells.Item(58, 31).Value = 0
$ws.Cells.Item(58, 32).Value = ""
$ws.Cells.Item(58, 33).Value = ""
$ws.Cells.Item(58, 34).Value = ""
$ws.Cells.Item(58, 35).Value = ""
$ws.Cells.Item(58, 36).Value = ""
$ws.Cells.Item(58, 37).Value = ""
$ws.Cells.Item(58, 38).Value = ""
$ws.Cells.Item(58, 39).Value = ""
$ws.Cells.Item(58, 40).Value = ""
$ws.Cells.Item(58, 41).Value = ""

# Row 59
$ws.Cells.Item(59, 1).Value = "WGE 97"
$ws.Cells.Item(59, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(59, 3).Value = "13-01-2026"
$ws.Cells.Item(59, 4).Value = 286962
$ws.Cells.Item(59, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(59, 6).Value = 34413429360
$ws.Cells.Item(59, 7).Value = "NEFT"
$ws.Cells.Item(59, 8).Value = "SBIN0003229"
$ws.Cells.Item(59, 9).Value = "AAAFW8862C"
$ws.Cells.Item(59, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(59, 11).Value = "Zillahmol S"
$ws.Cells.Item(59, 12).Value = "6c59549c-2fa0-41f7-a130-e04bedab34dd"
$ws.Cells.Item(59, 13).Value = 395502010022654
$ws.Cells.Item(59, 14).Value = "UBIN0826308"
$ws.Cells.Item(59, 15).Value = ""
$ws.Cells.Item(59, 16).Value = ""
$ws.Cells.Item(59, 17).Value = ""
$ws.Cells.Item(59, 18).Value = ""
$ws.Cells.Item(59, 19).Value = ""
$ws.Cells.Item(59, 20).Value = ""
$ws.Cells.Item(59, 21).Value = "pending"
$ws.Cells.Item(59, 22).Value = 1060.82
$ws.Cells.Item(59, 23).Value = ""
$ws.Cells.Item(59, 24).Value = "Hisham sir flat wifi (Today) RPA_UNIQUE_ID : 6bc0616f-1afe-499a-876f-aa605b4cecfb"
$ws.Cells.Item(59, 25).Value = "HO"
$ws.Cells.Item(59, 26).Value = "PAYMENT"
$ws.Cells.Item(59, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(59, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(59, 29).Value = 0
$ws.Cells.Item(59, 30).Value = 0
$ws.Cells.Item(59, 31).Value = 0
$ws.Cells.Item(59, 32).Value = ""
$ws.Cells.Item(59, 33).Value = ""
$ws.Cells.Item(59, 34).Value = ""
$ws.Cells.Item(59, 35).Value = ""
$ws.Cells.Item(59, 36).Value = ""
$ws.Cells.Item(59, 37).Value = ""
$ws.Cells.Item(59, 38).Value = ""
$ws.Cells.Item(59, 39).Value = ""
$ws.Cells.Item(59, 40).Value = ""
$ws.Cells.Item(59, 41).Value = ""

# Row 60
$ws.Cells.Item(60, 1).Value = "WGE 77"
$ws.Cells.Item(60, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(60, 3).Value = "13-01-2026"
$ws.Cells.Item(60, 4).Value = 286962
$ws.Cells.Item(60, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(60, 6).Value = 34413429360
$ws.Cells.Item(60, 7).Value = "NEFT"
$ws.Cells.Item(60, 8).Value = "SBIN0003229"
$ws.Cells.Item(60, 9).Value = "AAAFW8862C"
$ws.Cells.Item(60, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(60, 11).Value = "Anju M S"
$ws.Cells.Item(60, 12).Value = "5497c7c7-ae05-485b-9424-8251cdda7bc4"
$ws.Cells.Item(60, 13).Value = 345002010013320
$ws.Cells.Item(60, 14).Value = "UBIN0534501"
$ws.Cells.Item(60, 15).Value = ""
$ws.Cells.Item(60, 16).Value = ""
$ws.Cells.Item(60, 17).Value = ""
$ws.Cells.Item(60, 18).Value = ""
$ws.Cells.Item(60, 19).Value = ""
$ws.Cells.Item(60, 20).Value = ""
$ws.Cells.Item(60, 21).Value = "pending"
$ws.Cells.Item(60, 22).Value = 4300
$ws.Cells.Item(60, 23).Value = ""
$ws.Cells.Item(60, 24).Value = "Mumbai material purchase RPA_UNIQUE_ID : 84589d52-88fa-40d5-9bdc-1689163cf6c7"
$ws.Cells.Item(60, 25).Value = "Mumbai"
$ws.Cells.Item(60, 26).Value = "PAYMENT"
$ws.Cells.Item(60, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(60, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(60, 29).Value = 0
$ws.Cells.Item(60, 30).Value = 0
$ws.Cells.Item(60, 31).Value = 0
$ws.Cells.Item(60, 32).Value = ""
$ws.Cells.Item(60, 33).Value = ""
$ws.Cells.Item(60, 34).Value = ""
$ws.Cells.Item(60, 35).Value = ""
$ws.Cells.Item(60, 36).Value = ""
$ws.Cells.Item(60, 37).Value = ""
$ws.Cells.Item(60, 38).Value = ""
$ws.Cells.Item(60, 39).Value = ""
$ws.Cells.Item(60, 40).Value = ""
$ws.Cells.Item(60, 41).Value = ""

# Row 61
$ws.Cells.Item(61, 1).Value = "WGA019"
$ws.Cells.Item(61, 2).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(61, 3).Value = "13-01-2026"
$ws.Cells.Item(61, 4).Value = 286962
$ws.Cells.Item(61, 5).Value = "Western Interior Designers & Marine Contractors"
$ws.Cells.Item(61, 6).Value = 34413429360
$ws.Cells.Item(61, 7).Value = "NEFT"
$ws.Cells.Item(61, 8).Value = "SBIN0003229"
$ws.Cells.Item(61, 9).Value = "AAAFW8862C"
$ws.Cells.Item(61, 10).Value = "32AAAFW8862C1Z9"
$ws.Cells.Item(61, 11).Value = "TRAVEL DESIGNERS HUB"
$ws.Cells.Item(61, 12).Value = "6b806b49-1903-4ff5-891e-baf4234d264d"
$ws.Cells.Item(61, 13).Value = 10030200033897
$ws.Cells.Item(61, 14).Value = "FDRL0001003"
$ws.Cells.Item(61, 15).Value = ""
$ws.Cells.Item(61, 16).Value = ""
$ws.Cells.Item(61, 17).Value = ""
$ws.Cells.Item(61, 18).Value = ""
$ws.Cells.Item(61, 19).Value = ""
$ws.Cells.Item(61, 20).Value = ""
$ws.Cells.Item(61, 21).Value = "pending"
$ws.Cells.Item(61, 22).Value = 50000
$ws.Cells.Item(61, 23).Value = ""
$ws.Cells.Item(61, 24).Value = "Flight booking payments RPA_UNIQUE_ID : c9ac5b53-f187-456d-b876-0bc092c0aefd"
$ws.Cells.Item(61, 25).Value = "HO"
$ws.Cells.Item(61, 26).Value = "PAYMENT"
$ws.Cells.Item(61, 27).Value = "officeadmin@westernidc.com"
$ws.Cells.Item(61, 28).Value = "ESTIMATION NOT MATCHED"
$ws.Cells.Item(61, 29).Value = 0
$ws.Cells.Item(61, 30).Value = 0
$ws.Cells.Item(61, 31).Value = 0
$ws.Cells.Item(61, 32).Value = ""
$ws.Cells.Item(61, 33).Value = ""
$ws.Cells.Item(61, 34).Value = ""
$ws.Cells.Item(61, 35).Value = ""
$ws.Cells.Item(61, 36).Value = ""
$ws.Cells.Item(61, 37).Value = ""
$ws.Cells.Item(61, 38).Value = ""
$ws.Cells.Item(61, 39).Value = ""
$ws.Cells.Item(61, 40).Value = ""
$ws.Cells.Item(61, 41).Value = ""
